$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.99999998966469006
$ws.Range("A2").Value = 0.99635382945445561
$ws.Range("A3").Value = 0.98674008581422734
$ws.Range("A4").Value = 0.98768598773837446
$ws.Range("A5").Value = 0.97541881069704139
$ws.Range("A6").Value = 0.9453614446782832
$ws.Range("A7").Value = 0.94116119299592826
$ws.Range("A8").Value = 0.93536961297850363
$ws.Range("A9").Value = 0.92948405104361786
$ws.Range("A10").Value = 0.92451604391228726
$ws.Range("A11").Value = 0.923822116705733
$ws.Range("A12").Value = 0.92275744181873309
$ws.Range("A13").Value = 0.91991772339907363
$ws.Range("A14").Value = 0.91959745359927281
$ws.Range("A15").Value = 0.92054070368029817
$ws.Range("A16").Value = 0.92218483194887291
$ws.Range("A17").Value = 0.91847712834119077
$ws.Range("A18").Value = 0.91736824674795203
$ws.Range("A19").Value = 0.99426264282895538
$ws.Range("A20").Value = 0.98714580360245574
$ws.Range("A21").Value = 0.98574735227664356
$ws.Range("A22").Value = 0.98448285556292348
$ws.Range("A23").Value = 0.97970553253046844
$ws.Range("A24").Value = 0.96668489222070453
$ws.Range("A25").Value = 0.96022794878849083
$ws.Range("A26").Value = 0.94170651210704892
$ws.Range("A27").Value = 0.93686059494438312
$ws.Range("A28").Value = 0.91538984900019749
$ws.Range("A29").Value = 0.90011944098379648
$ws.Range("A30").Value = 0.89354931268155036
$ws.Range("A31").Value = 0.88589589831695292
$ws.Range("A32").Value = 0.88421664125192989
$ws.Range("A33").Value = 0.88369665261062624
